$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:D19")
$rng.Sort($ws.Range("A2:A19"), 1)
